# Gantt Chart.xlsx - "Add files via upload" edit
# Re-applies the task-tracker updates: bump the displayed week, refresh a
# few task progress %/end-dates, and fill in "TBD" assignees with real
# names now that they're known.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")

# --- Display Week: scroll the Gantt view forward two weeks (4 -> 6) -------
$ws.Range("E4").Value = 6

# --- Phase "Interface Design" (row 18): progress bumped 60% -> 80% --------
$ws.Range("D18").Value = 0.8

# --- Row 20 "Registration/log-in page": Kass joins John, progress + end ---
$ws.Range("C20").Value = "John, Kass"
$ws.Range("D20").Value = 0.85
$ws.Range("F20").Value = 43665

# --- Row 21 "Post Editing page": progress + end date pushed out -----------
$ws.Range("D21").Value = 0.85
$ws.Range("F21").Value = 43665

# --- Row 22 "DB Manipulation to Search Result": assignee + progress -------
$ws.Range("C22").Value = "James"
$ws.Range("D22").Value = 0.5

# --- Row 23 "Testing": assignee, progress, end date ------------------------
$ws.Range("C23").Value = "All"
$ws.Range("D23").Value = 0.3
$ws.Range("F23").Value = 43669

# --- Row 25 "Administrative Matter": assignee now known -------------------
$ws.Range("C25").Value = "John"

# --- Row 26 '"About Us" page': assignee now known --------------------------
$ws.Range("C26").Value = "John"

# --- Leave the cursor where the author last left it on save ---------------
[void]$ws.Range("J20").Select()
